# TC27_Canine_StudyUBC-AllBreeds_StageOfDisease.xlsx
# "Timing issue fix - keywords, updated tc1,2 in ubc01"
#
# The CasesTab query stored in B2 had an extraneous trailing `Cohort`
# output column that the other tab queries don't have (and that isn't
# matched anywhere useful) - drop that last RETURN line so TC1/TC2 line
# up with the rest of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in [ 'T3N0M1', 'Not Applicable']  OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Strip the trailing newline the here-string leaves behind so the value
# ends right after "Response to Treatment", matching the other queries'
# untrailed style.
$newCasesQuery = $newCasesQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newCasesQuery

# The row shrank by one wrapped line (the removed Cohort column), so the
# row height comes down from 304.5 to 290 - same as rows 3 and 4.
$ws.Rows.Item(2).RowHeight = 290

# Reset the sheet's saved view back to the top of the tab list and select
# the cell that was just edited.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
